$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting all the existing data down by one.
$ws.Rows(1).Insert()

# Populate the new row 1 with its two values.
$ws.Range("A1").Value = 26
$ws.Range("B1").Value = 52

# Move the selection to I5 (the old topLeftCell="R1"/selection Y30 pin is
# dropped along with it).
[void]$ws.Range("I5").Select()
